$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 539011121
$ws.Range("B3").Value = "Anna_Safonova_life"
$ws.Range("C3").Value = "3-минутная дыхательная пауза"
$ws.Range("D3").Value = "Breathing Space"
$ws.Range("E3").Value = "2025-10-19 12:34:22"
$ws.Range("F3").Value = "Получилось расслабиться"
$ws.Range("G3").Value = "Фокус на внутреннем мире"
$ws.Range("H3").Value = "Замедление"
$ws.Range("I3").Value = "Сложно успокоить мысли"
$ws.Range("J3").Value = "2025-10-19 12:34:22"

$wb.Save()
